# Insert two new rows of data at rows 23-24 (shifts all following rows down by 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("23:24").Insert()

# Row 23 (Calidad: Primera)
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(23, 3).Value = "Los Lagos"
$ws.Cells.Item(23, 4).Value = 44904
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100103
$ws.Cells.Item(23, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(23, 9).Value = 100103003
$ws.Cells.Item(23, 10).Value = "Damasco"
$ws.Cells.Item(23, 11).Value = "Castle Brite"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 600
$ws.Cells.Item(23, 14).Value = 21000
$ws.Cells.Item(23, 15).Value = 22000
$ws.Cells.Item(23, 16).Value = 21500
$ws.Cells.Item(23, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(23, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(23, 19).Value = 1344
$ws.Cells.Item(23, 20).Value = 16

# Row 24 (Calidad: Segunda)
$ws.Cells.Item(24, 1).Value = 4
$ws.Cells.Item(24, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(24, 3).Value = "Los Lagos"
$ws.Cells.Item(24, 4).Value = 44904
$ws.Cells.Item(24, 5).Value = 10
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100103
$ws.Cells.Item(24, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(24, 9).Value = 100103003
$ws.Cells.Item(24, 10).Value = "Damasco"
$ws.Cells.Item(24, 11).Value = "Castle Brite"
$ws.Cells.Item(24, 12).Value = "Segunda"
$ws.Cells.Item(24, 13).Value = 300
$ws.Cells.Item(24, 14).Value = 18000
$ws.Cells.Item(24, 15).Value = 18000
$ws.Cells.Item(24, 16).Value = 18000
$ws.Cells.Item(24, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(24, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(24, 19).Value = 1125
$ws.Cells.Item(24, 20).Value = 16

